$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 400; this shifts existing rows 400..494 down to 401..495
# and Excel carries the formatting (incl. the date style on column D) down from the row above.
$ws.Rows.Item(400).Insert()

# Populate the new row 400 with the new record's data.
# Columns A,B,C,E,F,G,H,I,N,Q,R repeat the same constant values used throughout this sheet;
# D (fecha), J (volumen), K/L/M (precios), O (origen) and P (precio $/Kg) are the new values.
$ws.Cells.Item(400, 1).Value = 5
$ws.Cells.Item(400, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(400, 3).Value = 'Maule'
$ws.Cells.Item(400, 4).Value = 45173
$ws.Cells.Item(400, 5).Value = 7
$ws.Cells.Item(400, 6).Value = 100112009
$ws.Cells.Item(400, 7).Value = 'Acelga'
$ws.Cells.Item(400, 8).Value = 'Sin especificar'
$ws.Cells.Item(400, 9).Value = 'Primera'
$ws.Cells.Item(400, 10).Value = 500
$ws.Cells.Item(400, 11).Value = 2000
$ws.Cells.Item(400, 12).Value = 2000
$ws.Cells.Item(400, 13).Value = 2000
$ws.Cells.Item(400, 14).Value = '$/docena de atados (4 kilos)'
$ws.Cells.Item(400, 15).Value = 'Provincia de Curicó'
$ws.Cells.Item(400, 16).Value = 500
$ws.Cells.Item(400, 17).Value = 4
$ws.Cells.Item(400, 18).Value = 'Hortaliza'
